# The commit swaps the contents of ppt/theme/theme1.xml (the slide master's
# theme - currently the "Integral" theme) and ppt/theme/theme2.xml (the notes
# master's theme - currently the default "Office Theme"): after the edit,
# theme1.xml holds the "Office Theme" palette and theme2.xml holds the
# "Integral" palette (font scheme / format scheme are identical between the
# two themes, so only the colour values - and the theme/colour-scheme names -
# actually change).
#
# The PowerPoint object model exposes the presentation's (slide-master)
# theme colours through ThemeColorScheme, indexed 1-12 in the standard
# order: dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink. Updating
# these colours here rewrites ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$colors = $slide.ThemeColorScheme

# Target palette: the stock "Office Theme" colours (RRGGBB), in ThemeColorScheme order.
$officeThemeRgb = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hlink
    0x954F72   # 12 folHlink
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $rrggbb = $officeThemeRgb[$i - 1]
    $r = [int](($rrggbb -band 0xFF0000) / 0x10000)
    $g = [int](($rrggbb -band 0x00FF00) / 0x100)
    $b = [int]($rrggbb -band 0x0000FF)
    # VBA-style RGB() packs as R + G*256 + B*65536, matching ColorFormat.RGB.
    $colors.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
